$d = $word.ActiveDocument

# Locate the "LOQ4084: ..." requirement paragraph. The three paragraphs that
# immediately follow it (a blank paragraph, the "Ver no Jupiter ..." line and
# the "(c) 2020 ..." footer line) are the ones being removed by this edit.
$anchorIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*LOQ4084*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -ne $null) {
    $firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
    $lastToRemove = $d.Paragraphs.Item($anchorIndex + 3)

    $removeRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $removeRange.Delete()
}
